$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet/title name to reflect new "through" date
$wb.Worksheets.Item(1).Name = "Through 2021-12-03"

# Update the December row label
$ws.Range("A13").Value = "December (through 12-03)"

# Update December row (row 13) values
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 11
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 9
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 13
$ws.Range("H13").Value = 20

# Update Total row (row 14) values
$ws.Range("B14").Value = 292
$ws.Range("C14").Value = 574
$ws.Range("D14").Value = 831
$ws.Range("E14").Value = 691
$ws.Range("F14").Value = 536
$ws.Range("G14").Value = 1277
$ws.Range("H14").Value = 1664
